$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Random Keys" column
$ws.Range("G1").Value = "Random Keys"

# Random key values (numbers) for rows 2-13
$ws.Range("G2").Value = 101101111
$ws.Range("G3").Value = 100110001
$ws.Range("G4").Value = 100000
$ws.Range("G5").Value = 111101011
$ws.Range("G6").Value = 101001111
$ws.Range("G7").Value = 101000010
$ws.Range("G8").Value = 1011111
$ws.Range("G9").Value = 1100100
$ws.Range("G10").Value = 101100010
$ws.Range("G11").Value = 11001101
$ws.Range("G12").Value = 110010001
$ws.Range("G13").Value = 10010101

# New row 14 - correct key label + value
$ws.Range("E14").Value = "correct key"
$ws.Range("F14").Value = 110100100

# Match final selection state from the authored workbook
$ws.Range("G13").Select()
